# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"   = @{ 2 = 733; 4 = 247; 5 = 2836; 6 = 56; 7 = 3756; 8 = 472; 9 = 947; 10 = 18 }
    "全部类型" = @{ 2 = 733; 5 = 247; 6 = 2836; 7 = 56; 8 = 3756; 9 = 472; 10 = 947; 11 = 18 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}

$wb.Save()
